# Insert 3 new weekly price rows (Chirimoya, Comercializadora del Agro de Limari)
# before the existing row 63, shifting the old rows 63:96 down to 66:99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63-96 down by 3 (inserts 3 blank rows at 63,64,65)
$ws.Rows("63:65").Insert()

# --- New row 63: Especial, $/bandeja 10 kilos ---
$ws.Range("A63").Value = 2
$ws.Range("B63").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 44825
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100107
$ws.Range("H63").Value = "Otros"
$ws.Range("I63").Value = 100107002
$ws.Range("J63").Value = "Chirimoya"
$ws.Range("K63").Value = "Cultivar IV Región"
$ws.Range("L63").Value = "Especial"
$ws.Range("M63").Value = 300
$ws.Range("N63").Value = 22000
$ws.Range("O63").Value = 23000
$ws.Range("P63").Value = 22500
$ws.Range("Q63").Value = "$/bandeja 10 kilos"
$ws.Range("R63").Value = "Provincia de Limarí"
$ws.Range("S63").Value = 2250
$ws.Range("T63").Value = 10

# --- New row 64: Primera, $/bandeja 10 kilos ---
$ws.Range("A64").Value = 2
$ws.Range("B64").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44825
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100107
$ws.Range("H64").Value = "Otros"
$ws.Range("I64").Value = 100107002
$ws.Range("J64").Value = "Chirimoya"
$ws.Range("K64").Value = "Cultivar IV Región"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 18000
$ws.Range("O64").Value = 19000
$ws.Range("P64").Value = 18500
$ws.Range("Q64").Value = "$/bandeja 10 kilos"
$ws.Range("R64").Value = "Provincia de Limarí"
$ws.Range("S64").Value = 1850
$ws.Range("T64").Value = 10

# --- New row 65: Segunda, $/bandeja 10 kilos ---
$ws.Range("A65").Value = 2
$ws.Range("B65").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44825
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100107
$ws.Range("H65").Value = "Otros"
$ws.Range("I65").Value = 100107002
$ws.Range("J65").Value = "Chirimoya"
$ws.Range("K65").Value = "Cultivar IV Región"
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 15000
$ws.Range("O65").Value = 16000
$ws.Range("P65").Value = 15500
$ws.Range("Q65").Value = "$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia de Limarí"
$ws.Range("S65").Value = 1550
$ws.Range("T65").Value = 10
